$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Shift the existing UML table from A1:D11 down/right to B2:E12 ---
$ws.Rows.Item(1).Insert()
$ws.Columns.Item(1).Insert()
$ws.Columns.Item(1).ColumnWidth = 4.75

# --- Add three new UML class boxes (Address, City, Country) starting at row 17 ---

# Row 17 (box headers / class names)
$ws.Range("B17").Value = "Address"
$ws.Range("C17").Value = "City"
$ws.Range("D17").Value = "Country"

# Row 18
$ws.Range("B18").Value = "addressId: int"
$ws.Range("C18").Value = "cityId: int"
$ws.Range("D18").Value = "countryId: int"

# Row 19
$ws.Range("B19").Value = "address: String"
$ws.Range("C19").Value = "city: String"
$ws.Range("D19").Value = "country: String"

# Row 20
$ws.Range("B20").Value = "address2: String"
$ws.Range("C20").Value = "countryId: int"
$ws.Range("D20").Value = "createDate: dateTime"

# Row 21
$ws.Range("B21").Value = "cityID: int"
$ws.Range("C21").Value = "createDate: dateTime"
$ws.Range("D21").Value = "createdBy: String"

# Row 22
$ws.Range("B22").Value = "postalCode: String"
$ws.Range("C22").Value = "createdBy: String"
$ws.Range("D22").Value = "lastUpdate: dateTime"

# Row 23
$ws.Range("B23").Value = "phone: String"
$ws.Range("C23").Value = "lastUpdate: dateTime"
$ws.Range("D23").Value = "lastUpdateBy: String"

# Row 24
$ws.Range("B24").Value = "createDate: dateTime"
$ws.Range("C24").Value = "lastUpdateBy: String"
$ws.Range("D24").Value = "set*()"

# Row 25
$ws.Range("B25").Value = "createdBy: String"
$ws.Range("C25").Value = "set*()"
$ws.Range("D25").Value = "get*()"

# Row 26
$ws.Range("B26").Value = "lastUpdate: dateTime"
$ws.Range("C26").Value = "get*()"

# Row 27
$ws.Range("B27").Value = "lastUpdateBy: String"

# Row 28
$ws.Range("B28").Value = "set*()"

# Row 29
$ws.Range("B29").Value = "get*()"

# --- Apply matching cell formatting by copying formats from the existing styled table ---
# style s=3 (box title, full border) source: C2
$ws.Range("C2").Copy()
$ws.Range("B17:D17").PasteSpecial(-4122)

# style s=4 (first attribute row under title) source: C3
$ws.Range("C3").Copy()
$ws.Range("B18:D18").PasteSpecial(-4122)

# style s=5 (attribute rows) source: C4
$ws.Range("C4").Copy()
$ws.Range("B19:D23").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("B24:C24").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("B25:B27").PasteSpecial(-4122)

# style s=6 (method-section divider row) source: C9
$ws.Range("C9").Copy()
$ws.Range("D24").PasteSpecial(-4122)
$ws.Range("C9").Copy()
$ws.Range("C25").PasteSpecial(-4122)
$ws.Range("C9").Copy()
$ws.Range("B28").PasteSpecial(-4122)

# style s=8 (box bottom row) source: C12
$ws.Range("C12").Copy()
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("C12").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("C12").Copy()
$ws.Range("B29").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Row-height touches to mirror the thick-bottom-border rows (15pt) ---
$ws.Rows.Item(16).RowHeight = 15
$ws.Rows.Item(25).RowHeight = 15
$ws.Rows.Item(26).RowHeight = 15
$ws.Rows.Item(29).RowHeight = 15

# --- Selection / active cell to match the final authored state ---
$ws.Range("E28").Select()
